$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-based values for rows 2-7 (columns A-T)
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf22"
$ws.Range("C2").Value = "Fgfr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.028883
$ws.Range("H2").Value = 0.086649
$ws.Range("I2").Value = 0.2501761214025038
$ws.Range("J2").Value = 0.2501761214025038
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.845768666666667
$ws.Range("N2").Value = 5.537306
$ws.Range("O2").Value = 0.01459089321241885
$ws.Range("P2").Value = 0.01459089321241885
$ws.Range("Q2").Value = 0.05331133639933334
$ws.Range("R2").Value = 0.479802027594
$ws.Range("S2").Value = 0.003650293071681068
$ws.Range("T2").Value = 0.003650293071681068

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf22"
$ws.Range("C3").Value = "Fgfr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.028883
$ws.Range("H3").Value = 0.086649
$ws.Range("I3").Value = 0.2501761214025038
$ws.Range("J3").Value = 0.2501761214025038
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 82.95722966666666
$ws.Range("N3").Value = 248.871689
$ws.Range("O3").Value = 0.6557810310272387
$ws.Range("P3").Value = 0.6557810310272387
$ws.Range("Q3").Value = 2.396053664462333
$ws.Range("R3").Value = 21.564482980161
$ws.Range("S3").Value = 0.1640607548317296
$ws.Range("T3").Value = 0.1640607548317296

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf22"
$ws.Range("C4").Value = "Fgfr1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.028883
$ws.Range("H4").Value = 0.086649
$ws.Range("I4").Value = 0.2501761214025038
$ws.Range("J4").Value = 0.2501761214025038
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 41.69841866666667
$ws.Range("N4").Value = 125.095256
$ws.Range("O4").Value = 0.3296280757603424
$ws.Range("P4").Value = 0.3296280757603424
$ws.Range("Q4").Value = 1.204375426349334
$ws.Range("R4").Value = 10.839378837144
$ws.Range("S4").Value = 0.08246507349909316
$ws.Range("T4").Value = 0.08246507349909316

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf22"
$ws.Range("C5").Value = "Fgfr1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.08656766666666667
$ws.Range("H5").Value = 0.259703
$ws.Range("I5").Value = 0.7498238785974961
$ws.Range("J5").Value = 0.7498238785974962
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.845768666666667
$ws.Range("N5").Value = 5.537306
$ws.Range("O5").Value = 0.01459089321241885
$ws.Range("P5").Value = 0.01459089321241885
$ws.Range("Q5").Value = 0.1597838866797778
$ws.Range("R5").Value = 1.438054980118
$ws.Range("S5").Value = 0.01094060014073778
$ws.Range("T5").Value = 0.01094060014073778

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf22"
$ws.Range("C6").Value = "Fgfr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.08656766666666667
$ws.Range("H6").Value = 0.259703
$ws.Range("I6").Value = 0.7498238785974961
$ws.Range("J6").Value = 0.7498238785974962
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 82.95722966666666
$ws.Range("N6").Value = 248.871689
$ws.Range("O6").Value = 0.6557810310272387
$ws.Range("P6").Value = 0.6557810310272387
$ws.Range("Q6").Value = 7.181413805374111
$ws.Range("R6").Value = 64.63272424836701
$ws.Range("S6").Value = 0.4917202761955091
$ws.Range("T6").Value = 0.4917202761955092

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf22"
$ws.Range("C7").Value = "Fgfr1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.08656766666666667
$ws.Range("H7").Value = 0.259703
$ws.Range("I7").Value = 0.7498238785974961
$ws.Range("J7").Value = 0.7498238785974962
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 39.34793966666667
$ws.Range("N7").Value = 118.043819
$ws.Range("O7").Value = 0.3296280757603424
$ws.Range("P7").Value = 0.3296280757603424
$ws.Range("Q7").Value = 3.609734807663111
$ws.Range("R7").Value = 32.487613268968
$ws.Range("S7").Value = 0.2471630022612492
$ws.Range("T7").Value = 0.2471630022612493
